$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 1857.0264
$ws.Range("J17").Value = 1873.081
$ws.Range("L17").Value = 5619.242999999999
$ws.Range("N17").Value = -5955.242999999999

# Row 28: The Writing Is Not on the Wall
$ws.Range("H28").Value = 847.9091
$ws.Range("I28").Value = 503.70587
$ws.Range("K28").Value = 503.70587
$ws.Range("M28").Value = -18.70587

# Row 74: Adhesive of Antipathy
$ws.Range("H74").Value = 6732.375
$ws.Range("I74").Value = 4057
$ws.Range("K74").Value = 4057
$ws.Range("M74").Value = -3121

# Row 77: It's Gonna Grow Back (L)
$ws.Range("H77").Value = 6732.375
$ws.Range("I77").Value = 4057
$ws.Range("K77").Value = 20285
$ws.Range("M77").Value = -15605

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 2497.818
$ws.Range("I132").Value = 2444.6843
$ws.Range("K132").Value = 7334.0529
$ws.Range("M132").Value = -4804.0529

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 2914.389
$ws.Range("I137").Value = 1907.6666
$ws.Range("J137").Value = 3417.75
$ws.Range("K137").Value = 5722.9998
$ws.Range("L137").Value = 10253.25
$ws.Range("M137").Value = -3172.9998
$ws.Range("N137").Value = -15353.25

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 2733.5881
$ws.Range("I2").Value = 3398.6365
$ws.Range("J2").Value = 1514.3334
$ws.Range("K2").Value = 3398.6365
$ws.Range("L2").Value = 1514.3334
$ws.Range("M2").Value = -3285.6365
$ws.Range("N2").Value = -1740.3334

# Row 98: Greaving
$ws.Range("H98").Value = 71163
$ws.Range("J98").Value = 71163
$ws.Range("L98").Value = 71163
$ws.Range("N98").Value = -77153

# Row 116: No Scope
$ws.Range("H116").Value = 2733.5881
$ws.Range("I116").Value = 3398.6365
$ws.Range("J116").Value = 1514.3334
$ws.Range("K116").Value = 3398.6365
$ws.Range("L116").Value = 1514.3334
$ws.Range("M116").Value = -1104.6365
$ws.Range("N116").Value = -6102.3334

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 4126.4287
$ws.Range("I122").Value = 3596.75
$ws.Range("K122").Value = 10790.25
$ws.Range("M122").Value = -8340.25

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 2733.5881
$ws.Range("I3").Value = 3398.6365
$ws.Range("J3").Value = 1514.3334
$ws.Range("K3").Value = 3398.6365
$ws.Range("L3").Value = 1514.3334
$ws.Range("M3").Value = -3284.6365
$ws.Range("N3").Value = -1742.3334

# Row 6: The Unkindest Cut
$ws.Range("H6").Value = 44680
$ws.Range("J6").Value = 44680
$ws.Range("L6").Value = 44680
$ws.Range("N6").Value = -44906

# Row 7: Thank You for Your Business
$ws.Range("H7").Value = 777.6667
$ws.Range("I7").Value = 99.8
$ws.Range("J7").Value = 1625
$ws.Range("K7").Value = 99.8
$ws.Range("L7").Value = 1625
$ws.Range("M7").Value = 13.2
$ws.Range("N7").Value = -1851

# Row 22: Riveting Run
$ws.Range("H22").Value = 483.82352
$ws.Range("I22").Value = 483.82352
$ws.Range("K22").Value = 483.82352
$ws.Range("M22").Value = -310.82352

# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 1326.1428
$ws.Range("I86").Value = 1266.6
$ws.Range("K86").Value = 1266.6
$ws.Range("M86").Value = -143.5999999999999

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 1326.1428
$ws.Range("I89").Value = 1266.6
$ws.Range("K89").Value = 6333
$ws.Range("M89").Value = -717

# Row 107: The Gold Experience
$ws.Range("H107").Value = 4953.2
$ws.Range("J107").Value = 4253.857
$ws.Range("L107").Value = 4253.857
$ws.Range("N107").Value = -8093.857

# Row 135: Axes to the Maxes
$ws.Range("H135").Value = 103322.336
$ws.Range("J135").Value = 103322.336
$ws.Range("L135").Value = 103322.336
$ws.Range("N135").Value = -113462.336

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 1499.8636
$ws.Range("J22").Value = 2702.7778
$ws.Range("L22").Value = 2702.7778
$ws.Range("N22").Value = -3402.7778

# Row 107: Built to Last
$ws.Range("H107").Value = 1850.1333
$ws.Range("I107").Value = 2130
$ws.Range("J107").Value = 730.6667
$ws.Range("K107").Value = 2130
$ws.Range("L107").Value = 730.6667
$ws.Range("M107").Value = -210
$ws.Range("N107").Value = -4570.6667

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 1874.1562
$ws.Range("I122").Value = 1678.5358
$ws.Range("J122").Value = 3243.5
$ws.Range("K122").Value = 5035.607400000001
$ws.Range("L122").Value = 9730.5
$ws.Range("M122").Value = -2585.607400000001
$ws.Range("N122").Value = -14630.5

# Row 123: A Real Grind
$ws.Range("H123").Value = 229999.2
$ws.Range("J123").Value = 229999.2
$ws.Range("L123").Value = 229999.2
$ws.Range("N123").Value = -239799.2

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2388.5264
$ws.Range("I132").Value = 1959.6875
$ws.Range("K132").Value = 5879.0625
$ws.Range("M132").Value = -3349.0625

$ws = $wb.Worksheets.Item("CUL")
# Row 56: Culture Club
$ws.Range("H56").Value = 7485.5713
$ws.Range("I56").Value = 7485.5713
$ws.Range("K56").Value = 7485.5713
$ws.Range("M56").Value = -6955.5713

# Row 112: Sweet Tooth
$ws.Range("H112").Value = 10000
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 10000
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 30000
$ws.Range("M112").Value = ""
$ws.Range("N112").Value = -32216

# Row 125: At Any Temperature
$ws.Range("H125").Value = 19232.4
$ws.Range("I125").Value = 2030
$ws.Range("J125").Value = 23533
$ws.Range("K125").Value = 6090
$ws.Range("L125").Value = 70599
$ws.Range("M125").Value = -1170
$ws.Range("N125").Value = -80439

# Row 134: Don't Knock It Till You've Tried It
$ws.Range("H134").Value = 3530
$ws.Range("I134").Value = 2287.5
$ws.Range("J134").Value = 8500
$ws.Range("K134").Value = 6862.5
$ws.Range("L134").Value = 25500
$ws.Range("M134").Value = -1792.5
$ws.Range("N134").Value = -35640

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 2453.1
$ws.Range("I102").Value = 1604.1333
$ws.Range("K102").Value = 1604.1333
$ws.Range("M102").Value = 17.86670000000004

# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 4572.7144
$ws.Range("I113").Value = 3499.5
$ws.Range("J113").Value = 5002
$ws.Range("K113").Value = 3499.5
$ws.Range("L113").Value = 5002
$ws.Range("M113").Value = -1329.5
$ws.Range("N113").Value = -9342

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 5335.2666
$ws.Range("I122").Value = 2754.8333
$ws.Range("K122").Value = 8264.499899999999
$ws.Range("M122").Value = -5814.499899999999

# Row 132: On Board for Lar
$ws.Range("H132").Value = 1647.4572
$ws.Range("I132").Value = 1624.4062
$ws.Range("K132").Value = 4873.2186
$ws.Range("M132").Value = -2343.2186

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 3693.3125
$ws.Range("I7").Value = 3554.889
$ws.Range("J7").Value = 3871.2856
$ws.Range("K7").Value = 3554.889
$ws.Range("L7").Value = 3871.2856
$ws.Range("M7").Value = -3442.889
$ws.Range("N7").Value = -4095.2856

# Row 40: Best Served Toad
$ws.Range("H40").Value = 5841.1333
$ws.Range("I40").Value = 6225.375
$ws.Range("J40").Value = 5402
$ws.Range("K40").Value = 6225.375
$ws.Range("L40").Value = 5402
$ws.Range("M40").Value = -6089.375
$ws.Range("N40").Value = -5674

# Row 122: Hell on Leather
$ws.Range("H122").Value = 3269.6316
$ws.Range("I122").Value = 2979.5625
$ws.Range("J122").Value = 4816.6665
$ws.Range("K122").Value = 8938.6875
$ws.Range("L122").Value = 14449.9995
$ws.Range("M122").Value = -6488.6875
$ws.Range("N122").Value = -19349.9995

# Row 126: Battered Books
$ws.Range("H126").Value = 3693.3125
$ws.Range("I126").Value = 3554.889
$ws.Range("J126").Value = 3871.2856
$ws.Range("K126").Value = 10664.667
$ws.Range("L126").Value = 11613.8568
$ws.Range("M126").Value = -8194.667000000001
$ws.Range("N126").Value = -16553.8568

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 3128.7576
$ws.Range("I132").Value = 2970.6897
$ws.Range("J132").Value = 4274.75
$ws.Range("K132").Value = 8912.069100000001
$ws.Range("L132").Value = 12824.25
$ws.Range("M132").Value = -6382.069100000001
$ws.Range("N132").Value = -17884.25

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 6532.1816
$ws.Range("J136").Value = 12500
$ws.Range("L136").Value = 37500
$ws.Range("N136").Value = -42600

$ws = $wb.Worksheets.Item("WVR")
# Row 42: Put on Your Party Pants
$ws.Range("H42").Value = 50000
$ws.Range("J42").Value = 50000
$ws.Range("L42").Value = 50000
$ws.Range("N42").Value = -50756

# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 1349.5
$ws.Range("I81").Value = 1349.5
$ws.Range("K81").Value = 2699
$ws.Range("M81").Value = -1638

# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 1349.5
$ws.Range("I84").Value = 1349.5
$ws.Range("K84").Value = 13495
$ws.Range("M84").Value = -8191

# Row 122: Heavy Armoire
$ws.Range("H122").Value = 2815.9167
$ws.Range("I122").Value = 2754.5557
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 8263.667099999999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -5813.667099999999
$ws.Range("N122").Value = -13900

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 3180.182
$ws.Range("I126").Value = 3180.182
$ws.Range("K126").Value = 9540.545999999998
$ws.Range("M126").Value = -7070.545999999998

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 4195.7427
$ws.Range("I132").Value = 2607.6
$ws.Range("K132").Value = 7822.799999999999
$ws.Range("M132").Value = -5292.799999999999
